# Working on integrating new wind incentives -- #495
#
# Insert two new columns ("min_size_kw" / "max_size_kw") between the
# existing "sector_abbr" and "val_pct_cost" columns, populate the new
# column values for the one data row, and nudge a couple of cosmetic
# view settings (active selection, workbook window width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at D:E -- everything from the old D onward
# (val_pct_cost, cap_dlrs, exp_date, dsire_program_name,
# dsire_last_updated, dsire_link) slides right by two columns.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("D1").Value = "min_size_kw"
$ws.Range("E1").Value = "max_size_kw"

# New data cells for the single data row.
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "Inf"

# Match column widths with the neighboring sector_abbr column.
$ws.Range("D1:E1").ColumnWidth = $ws.Range("C1").ColumnWidth

# Move the active selection (matches the post-edit cursor position
# recorded by Excel).
$ws.Range("D3").Select()

# Widen the saved workbook window.
$excel.ActiveWindow.Width = 38400
